$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item('Citywide Totals')
$ws.Range("K2").Value = 2528
$ws.Range("K3").Value = 2449
$ws.Range("J4").Value = 581
$ws.Range("K4").Value = 512
$ws.Range("K5").Value = 162
$ws.Range("K6").Value = 3048
$ws.Range("J7").Value = 8616
$ws.Range("K7").Value = 8699

$ws = $wb.Worksheets.Item('Logan Square')
$ws.Range("K6").Value = 67
$ws.Range("K7").Value = 128

$ws = $wb.Worksheets.Item('Austin')
$ws.Range("K2").Value = 169
$ws.Range("K3").Value = 174
$ws.Range("K6").Value = 193
$ws.Range("K7").Value = 580

$ws = $wb.Worksheets.Item('South Chicago')
$ws.Range("K6").Value = 44
$ws.Range("K7").Value = 192

$ws = $wb.Worksheets.Item('Garfield Park')
$ws.Range("K2").Value = 97
$ws.Range("K3").Value = 124
$ws.Range("K6").Value = 94
$ws.Range("K7").Value = 343

$ws = $wb.Worksheets.Item('Grand Crossing')
$ws.Range("K3").Value = 97
$ws.Range("K6").Value = 85
$ws.Range("K7").Value = 279

$ws = $wb.Worksheets.Item('New City')
$ws.Range("K6").Value = 83
$ws.Range("K7").Value = 205

$ws = $wb.Worksheets.Item('Woodlawn')
$ws.Range("K4").Value = 7
$ws.Range("K6").Value = 43
$ws.Range("K7").Value = 155

$ws = $wb.Worksheets.Item('Fuller Park')
$ws.Range("K6").Value = 9
$ws.Range("K7").Value = 25

$ws = $wb.Worksheets.Item('By Neighborhood')
$ws.Range("K2").Value = 62
$ws.Range("K4").Value = 33
$ws.Range("K6").Value = 71
$ws.Range("K7").Value = 257
$ws.Range("K8").Value = 580
$ws.Range("K14").Value = 51
$ws.Range("K15").Value = 87
$ws.Range("K16").Value = 26
$ws.Range("K18").Value = 58
$ws.Range("K19").Value = 255
$ws.Range("K20").Value = 198
$ws.Range("K21").Value = 23
$ws.Range("K23").Value = 78
$ws.Range("K29").Value = 441
$ws.Range("K30").Value = 25
$ws.Range("K31").Value = 99
$ws.Range("K33").Value = 343
$ws.Range("K37").Value = 279
$ws.Range("K42").Value = 305
$ws.Range("K44").Value = 84
$ws.Range("K47").Value = 48
$ws.Range("K48").Value = 109
$ws.Range("K51").Value = 95
$ws.Range("K52").Value = 237
$ws.Range("K53").Value = 128
$ws.Range("K54").Value = 162
$ws.Range("K60").Value = 57
$ws.Range("J63").Value = 46
$ws.Range("K63").Value = 34
$ws.Range("K64").Value = 56
$ws.Range("K65").Value = 205
$ws.Range("K66").Value = 30
$ws.Range("K67").Value = 341
$ws.Range("K70").Value = 16
$ws.Range("K72").Value = 41
$ws.Range("K74").Value = 11
$ws.Range("K76").Value = 126
$ws.Range("K78").Value = 119
$ws.Range("K83").Value = 192
$ws.Range("K84").Value = 63
$ws.Range("K85").Value = 418
$ws.Range("K86").Value = 55
$ws.Range("K89").Value = 114
$ws.Range("K91").Value = 83
$ws.Range("K99").Value = 155
$ws.Range("J101").Value = 8616
$ws.Range("K101").Value = 8699

$ws = $wb.Worksheets.Item('Gage Park')
$ws.Range("K6").Value = 39
$ws.Range("K7").Value = 99

$ws = $wb.Worksheets.Item('North Lawndale')
$ws.Range("K3").Value = 106
$ws.Range("K7").Value = 341

$ws = $wb.Worksheets.Item('South Deering')
$ws.Range("K3").Value = 20
$ws.Range("K7").Value = 63

$ws = $wb.Worksheets.Item('Loop')
$ws.Range("K4").Value = 9
$ws.Range("K7").Value = 162

$ws = $wb.Worksheets.Item('Englewood')
$ws.Range("K2").Value = 119
$ws.Range("K3").Value = 147
$ws.Range("K4").Value = 25
$ws.Range("K7").Value = 441

$ws = $wb.Worksheets.Item('Lake View')
$ws.Range("K3").Value = 20
$ws.Range("K7").Value = 109

$ws = $wb.Worksheets.Item('Chatham')
$ws.Range("K2").Value = 80
$ws.Range("K3").Value = 68
$ws.Range("K7").Value = 255

$ws = $wb.Worksheets.Item('Irving Park')
$ws.Range("K3").Value = 23
$ws.Range("K7").Value = 84

$ws = $wb.Worksheets.Item('River North')
$ws.Range("K6").Value = 77
$ws.Range("K7").Value = 126

$ws = $wb.Worksheets.Item('Bridgeport')
$ws.Range("K2").Value = 21
$ws.Range("K3").Value = 10
$ws.Range("K7").Value = 51

$ws = $wb.Worksheets.Item('Ashburn')
$ws.Range("K2").Value = 25
$ws.Range("K7").Value = 71

$ws = $wb.Worksheets.Item('Humboldt Park')
$ws.Range("K3").Value = 95
$ws.Range("K7").Value = 305

$ws = $wb.Worksheets.Item('Rogers Park')
$ws.Range("K3").Value = 28
$ws.Range("K6").Value = 45
$ws.Range("K7").Value = 119

$ws = $wb.Worksheets.Item('Douglas')
$ws.Range("K3").Value = 23
$ws.Range("K7").Value = 78

$ws = $wb.Worksheets.Item('Washington Park')
$ws.Range("K3").Value = 35
$ws.Range("K7").Value = 83

$ws = $wb.Worksheets.Item('Chinatown')
$ws.Range("K6").Value = 13
$ws.Range("K7").Value = 23

$ws = $wb.Worksheets.Item('Near South Side')
$ws.Range("K6").Value = 19
$ws.Range("K7").Value = 56

$ws = $wb.Worksheets.Item('Chicago Lawn')
$ws.Range("K3").Value = 56
$ws.Range("K6").Value = 68
$ws.Range("K7").Value = 198

$ws = $wb.Worksheets.Item('Calumet Heights')
$ws.Range("K3").Value = 19
$ws.Range("K5").Value = 2
$ws.Range("K7").Value = 58

$ws = $wb.Worksheets.Item('Auburn Gresham')
$ws.Range("K3").Value = 82
$ws.Range("K6").Value = 69
$ws.Range("K7").Value = 257

$ws = $wb.Worksheets.Item('Kenwood')
$ws.Range("K6").Value = 15
$ws.Range("K7").Value = 48

$ws = $wb.Worksheets.Item('Brighton Park')
$ws.Range("K2").Value = 30
$ws.Range("K7").Value = 87

$ws = $wb.Worksheets.Item('North Center')
$ws.Range("K3").Value = 5
$ws.Range("K6").Value = 30

$ws = $wb.Worksheets.Item('Albany Park')
$ws.Range("K3").Value = 17
$ws.Range("K6").Value = 23
$ws.Range("K7").Value = 62

$ws = $wb.Worksheets.Item('O''Hare')
$ws.Range("K2").Value = 6
$ws.Range("K6").Value = 16

$ws = $wb.Worksheets.Item('Uptown')
$ws.Range("K2").Value = 27
$ws.Range("K7").Value = 114

$ws = $wb.Worksheets.Item('Streeterville')
$ws.Range("K3").Value = 12
$ws.Range("K6").Value = 55

$ws = $wb.Worksheets.Item('Little Italy, UIC')
$ws.Range("K2").Value = 25
$ws.Range("K7").Value = 95

$ws = $wb.Worksheets.Item('Morgan Park')
$ws.Range("K2").Value = 16
$ws.Range("K7").Value = 57

$ws = $wb.Worksheets.Item('South Shore')
$ws.Range("K2").Value = 155
$ws.Range("K3").Value = 142
$ws.Range("K7").Value = 418

$ws = $wb.Worksheets.Item('Old Town')
$ws.Range("K5").Value = 24
$ws.Range("K6").Value = 41

$ws = $wb.Worksheets.Item('Little Village')
$ws.Range("K6").Value = 99
$ws.Range("K7").Value = 237

$ws = $wb.Worksheets.Item('Archer Heights')
$ws.Range("K3").Value = 10
$ws.Range("K7").Value = 33

$ws = $wb.Worksheets.Item('Bucktown')
$ws.Range("K2").Value = 6
$ws.Range("K7").Value = 26

$ws = $wb.Worksheets.Item('East Village')
$ws.Range("K5").Value = 8
$ws.Range("K6").Value = 11
